$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.49"
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Value = "'1.72%"
$ws.Range("E2").Style = $ws.Range("B2").Style
$ws.Range("D3").Value = "'27.15"
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Value = "'1.65%"
$ws.Range("E3").Style = $ws.Range("B3").Style
$ws.Range("D4").Value = "'4.903"
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").Value = "'0.06%"
$ws.Range("E4").Style = $ws.Range("B4").Style
$ws.Range("D5").Value = "'0.06417"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "'1.60%"
$ws.Range("E5").Style = $ws.Range("B5").Style
$ws.Range("D6").Value = "'6.949"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "'0.62%"
$ws.Range("E6").Style = $ws.Range("B6").Style
$ws.Range("D7").Value = "'1.241"
$ws.Range("D7").Style = $ws.Range("B7").Style
$ws.Range("E7").Value = "'-6.58%"
$ws.Range("E7").Style = $ws.Range("B7").Style
$ws.Range("D8").Value = "'0.8825"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "'-0.41%"
$ws.Range("E8").Style = $ws.Range("B8").Style
$ws.Range("D9").Value = "'0.1522"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "'4.19%"
$ws.Range("E9").Style = $ws.Range("B9").Style
$ws.Range("D10").Value = "'0.05027"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "'-2.53%"
$ws.Range("E10").Style = $ws.Range("B10").Style
$ws.Range("D11").Value = "'0.07503"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "'1.62%"
$ws.Range("E11").Style = $ws.Range("B11").Style
$ws.Range("D12").Value = "'0.02910"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "'-7.00%"
$ws.Range("E12").Style = $ws.Range("B12").Style
$ws.Range("D13").Value = "'0.09006"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "'-0.33%"
$ws.Range("E13").Style = $ws.Range("B13").Style
$ws.Range("D14").Value = "'0.001567"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "'-0.45%"
$ws.Range("E14").Style = $ws.Range("B14").Style
$ws.Range("D15").Value = "'0.0006412"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "'1.71%"
$ws.Range("E15").Style = $ws.Range("B15").Style
$ws.Range("D16").Value = "'0.005708"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "'-5.26%"
$ws.Range("E16").Style = $ws.Range("B16").Style
$ws.Range("D17").Value = "'3.460"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "'-0.06%"
$ws.Range("E17").Style = $ws.Range("B17").Style
$ws.Range("D18").Value = "'3.315"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "'-1.37%"
$ws.Range("E18").Style = $ws.Range("B18").Style
$ws.Range("E19").Value = "'0.03%"
$ws.Range("E19").Style = $ws.Range("B19").Style
$ws.Range("E20").Value = "'-0.95%"
$ws.Range("E20").Style = $ws.Range("B20").Style
$ws.Range("D21").Value = "'0.1337"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "'1.47%"
$ws.Range("E21").Style = $ws.Range("B21").Style
$ws.Range("D22").Value = "'3.913"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "'0.31%"
$ws.Range("E22").Style = $ws.Range("B22").Style
$ws.Range("E23").Value = "'1.76%"
$ws.Range("E23").Style = $ws.Range("B23").Style
$ws.Range("E24").Value = "'-0.32%"
$ws.Range("E24").Style = $ws.Range("B24").Style
$ws.Range("E25").Value = "'5.07%"
$ws.Range("E25").Style = $ws.Range("B25").Style
$ws.Range("E26").Value = "'0.04%"
$ws.Range("E26").Style = $ws.Range("B26").Style
$ws.Range("E27").Value = "'14.06%"
$ws.Range("E27").Style = $ws.Range("B27").Style
$ws.Range("D40").Value = "'0.04141"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "'2.94%"
$ws.Range("E40").Style = $ws.Range("B40").Style
$ws.Range("D41").Value = "'0.006807"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "'2.86%"
$ws.Range("E41").Style = $ws.Range("B41").Style
$ws.Range("D42").Value = "'0.1176"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E43").Value = "'13.92%"
$ws.Range("E43").Style = $ws.Range("B43").Style
$ws.Range("D44").Value = "'0.01171"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "'-4.39%"
$ws.Range("E44").Style = $ws.Range("B44").Style
$ws.Range("D45").Value = "'0.00005203"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "'-2.11%"
$ws.Range("E45").Style = $ws.Range("B45").Style
$ws.Range("D46").Value = "'1.486"
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = "'-36.89%"
$ws.Range("E46").Style = $ws.Range("B46").Style
$ws.Range("E47").Value = "'-22.14%"
$ws.Range("E47").Style = $ws.Range("B47").Style
